$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "measurement_year"

# Add header for new column F
$ws.Range("F1").Value = "microbes"

# Add the new microbes data values (F3:F6 = 6, F2 left empty)
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 6

# Update selection to F2 to match final cursor position
$ws.Range("F2").Select()
